$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new shared string / header in K1
$ws.Range("K1").Value = "Cleaning Data in Python"
$ws.Range("K1").Style = $ws.Range("J1").Style

# Add new data value in K2
$ws.Range("K2").Value = 3

# Set column K width to match J's bestFit-style formatting (col style index 2)
$ws.Range("K2").Style = $ws.Range("J2").Style
$ws.Columns("K").ColumnWidth = 22.42578125

# Update selection / view
$ws.Application.ActiveWindow.ScrollColumn = 3
$ws.Range("K3").Select()
